$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 510
$ws.Range("I12").Value = 510
$ws.Range("K12").Value = 510
$ws.Range("M12").Value = -340

# Row 32
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2674

# Row 33
$ws.Range("H33").Value = 176.5
$ws.Range("I33").Value = 176.5
$ws.Range("K33").Value = 176.5
$ws.Range("M33").Value = 52.5

# Row 106
$ws.Range("H106").Value = 43413.285
$ws.Range("I106").Value = 43413.285
$ws.Range("K106").Value = 43413.285
$ws.Range("M106").Value = -42782.285

# Row 138
$ws.Range("H138").Value = 5525.9736
$ws.Range("J138").Value = 4656.091
$ws.Range("L138").Value = 13968.273
$ws.Range("N138").Value = -24248.273


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7976.2393
$ws.Range("I32").Value = 5674.0586
$ws.Range("K32").Value = 5674.0586
$ws.Range("M32").Value = -5387.0586

# Row 45
$ws.Range("H45").Value = 1859
$ws.Range("I45").Value = 1719.6
$ws.Range("J45").Value = 2091.3333
$ws.Range("K45").Value = 1719.6
$ws.Range("L45").Value = 2091.3333
$ws.Range("M45").Value = -1342.6
$ws.Range("N45").Value = -2845.3333

# Row 61
$ws.Range("H61").Value = 2130.111
$ws.Range("I61").Value = 2016.4
$ws.Range("K61").Value = 2016.4
$ws.Range("M61").Value = -1804.4

# Row 122
$ws.Range("H122").Value = 4306.923
$ws.Range("I122").Value = 3998.889
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11996.667
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -9546.667000000001
$ws.Range("N122").Value = -19900

# Row 132
$ws.Range("H132").Value = 1509.0256
$ws.Range("I132").Value = 1418.875
$ws.Range("J132").Value = 1921.1428
$ws.Range("K132").Value = 4256.625
$ws.Range("L132").Value = 5763.428400000001
$ws.Range("M132").Value = -1726.625
$ws.Range("N132").Value = -10823.4284

# Row 136
$ws.Range("H136").Value = 2130.111
$ws.Range("I136").Value = 2016.4
$ws.Range("K136").Value = 6049.200000000001
$ws.Range("M136").Value = -3499.200000000001


$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1151.6
$ws.Range("I94").Value = 1168.4445
$ws.Range("K94").Value = 1168.4445
$ws.Range("M94").Value = -717.4445000000001

# Row 102
$ws.Range("H102").Value = 15000
$ws.Range("I102").Value = 15000
$ws.Range("K102").Value = 15000
$ws.Range("M102").Value = -11755

# Row 107
$ws.Range("H107").Value = 1333.3334
$ws.Range("I107").Value = 1350
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1350
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 570
$ws.Range("N107").Value = -5040

# Row 141
$ws.Range("H141").Value = 9999.5
$ws.Range("I141").Value = 9999.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9999.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4819.5
$ws.Range("N141").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2862.4138
$ws.Range("I58").Value = 1201.8334
$ws.Range("K58").Value = 1201.8334
$ws.Range("M58").Value = -998.8334

# Row 107
$ws.Range("H107").Value = 737
$ws.Range("I107").Value = 805.5
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 805.5
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1114.5
$ws.Range("N107").Value = -4440

# Row 122
$ws.Range("H122").Value = 2563.7827
$ws.Range("I122").Value = 3433.9285
$ws.Range("J122").Value = 1210.2222
$ws.Range("K122").Value = 10301.7855
$ws.Range("L122").Value = 3630.6666
$ws.Range("M122").Value = -7851.7855
$ws.Range("N122").Value = -8530.6666

# Row 132
$ws.Range("H132").Value = 2132.5527
$ws.Range("I132").Value = 1919.9459
$ws.Range("K132").Value = 5759.8377
$ws.Range("M132").Value = -3229.8377

# Row 134
$ws.Range("H134").Value = 2565
$ws.Range("I134").Value = 1446.1
$ws.Range("K134").Value = 4338.299999999999
$ws.Range("M134").Value = -1803.299999999999

# Row 136
$ws.Range("H136").Value = 2862.4138
$ws.Range("I136").Value = 1201.8334
$ws.Range("K136").Value = 3605.5002
$ws.Range("M136").Value = -1055.5002


$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 399.7143
$ws.Range("I92").Value = 452.69232
$ws.Range("K92").Value = 1358.07696
$ws.Range("M92").Value = -110.0769599999999

# Row 113
$ws.Range("H113").Value = 1121.3334
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1121.3334
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3364.0002
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -7704.0002

# Row 132
$ws.Range("H132").Value = 11996.4
$ws.Range("I132").Value = 11995.5
$ws.Range("K132").Value = 107959.5
$ws.Range("M132").Value = -105429.5


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3425
$ws.Range("J80").Value = 3500
$ws.Range("L80").Value = 3500
$ws.Range("N80").Value = -5496

# Row 83
$ws.Range("H83").Value = 3425
$ws.Range("J83").Value = 3500
$ws.Range("L83").Value = 17500
$ws.Range("N83").Value = -27484

# Row 107
$ws.Range("H107").Value = 581.5
$ws.Range("I107").Value = 581.5
$ws.Range("K107").Value = 581.5
$ws.Range("M107").Value = 1338.5


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2306.6667
$ws.Range("I7").Value = 2306.6667
$ws.Range("K7").Value = 2306.6667
$ws.Range("M7").Value = -2194.6667

# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 61
$ws.Range("H61").Value = 4679.3125
$ws.Range("I61").Value = 4679.3125
$ws.Range("K61").Value = 4679.3125
$ws.Range("M61").Value = -4477.3125

# Row 68
$ws.Range("H68").Value = 2410.375
$ws.Range("J68").Value = 1947
$ws.Range("L68").Value = 1947
$ws.Range("N68").Value = -3445

# Row 71
$ws.Range("H71").Value = 2410.375
$ws.Range("J71").Value = 1947
$ws.Range("L71").Value = 9735
$ws.Range("N71").Value = -17223

# Row 113
$ws.Range("H113").Value = 4679.3125
$ws.Range("I113").Value = 4679.3125
$ws.Range("K113").Value = 4679.3125
$ws.Range("M113").Value = -2509.3125

# Row 126
$ws.Range("H126").Value = 2306.6667
$ws.Range("I126").Value = 2306.6667
$ws.Range("K126").Value = 6920.000100000001
$ws.Range("M126").Value = -4450.000100000001

# Row 141
$ws.Range("H141").Value = 61428.25
$ws.Range("I141").Value = 48332.668
$ws.Range("J141").Value = 100715
$ws.Range("K141").Value = 48332.668
$ws.Range("L141").Value = 100715
$ws.Range("M141").Value = -43152.668
$ws.Range("N141").Value = -111075


$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16040

# Row 60
$ws.Range("H60").Value = 38996
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

# Row 81
$ws.Range("H81").Value = 10398.5
$ws.Range("I81").Value = 10398.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 20797
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -19736
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 10398.5
$ws.Range("I84").Value = 10398.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 103985
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -98681
$ws.Range("N84").ClearContents()

